$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows to append after the existing last row (238).
$data = @(
    @(44313, 0, 14, 90.42175289026675),
    @(44314, 1, 13, 83.96305625524769),
    @(44315, 3, 14, 90.42175289026675),
    @(44316, 9, 18, 116.256539430343),
    @(44317, 4, 20, 129.1739327003811),
    @(44318, 2, 19, 122.715236065362)
)

$startRow = 239
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the formatting of the last existing date cell (A238) onto the
    # new date cell so it keeps the same style (s="2": centered/bold/border
    # plus the date number format) before writing the value.
    $ws.Range("A238").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $row[0]

    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}
